$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "최종점수" (column K) values
$ws.Range("K2").Value = 62.7
$ws.Range("K3").Value = 54.5
$ws.Range("K4").Value = 50.7
$ws.Range("K5").Value = 47.9
$ws.Range("K6").Value = 39.9
$ws.Range("K7").Value = 39.9

# Update "MACRO_SCORE" (column N) values
$ws.Range("N2").Value = 51.53902399942638
$ws.Range("N3").Value = 51.53902399942638
$ws.Range("N4").Value = 51.53902399942638
$ws.Range("N5").Value = 51.53902399942638
$ws.Range("N6").Value = 51.53902399942638
$ws.Range("N7").Value = 51.53902399942638
